# Commit: "model needs to include both initial # of cells occupied and total
# initial biomass"
#
# The existing "initial01".."initial04" columns (N:Q) represented a single
# initial-condition number per species. The model now needs TWO numbers per
# species: the initial number of cells occupied ("initialNNcells") and the
# initial total biomass ("initialNNtotmass"). We keep the original
# initial01..04 columns (and their values) in place - renaming them to the
# "...totmass" variant - and insert four brand-new "...cells" columns in
# front of them, seeded with the same placeholder values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert 4 blank columns right where "initial01" used to live (column N).
#    Everything from N onward (maxrgr01.., overwinter01.., halfsat*, uptake*)
#    shifts right by 4 columns; the untouched former N:Q ("initial01"..
#    "initial04") ends up at R:U.
$ws.Range("N1:Q1").EntireColumn.Insert(-4161)

# Inserting columns makes Excel carry over the formatting of the column to
# the left (M, numbspecies, style index 2) onto the brand-new columns. The
# source data never styled these columns, so reset them back to the
# workbook's default "Normal" style.
$ws.Range("N1:Q7").Style = "Normal"

# 2) New N1:Q1 headers - the "# of cells occupied" columns.
$ws.Range("N1").Value = "initial01cells"
$ws.Range("O1").Value = "initial02cells"
$ws.Range("P1").Value = "initial03cells"
$ws.Range("Q1").Value = "initial04cells"

# 3) Seed the new columns with the same placeholder values used by the
#    original "initial01".."initial04" columns (now shifted to R:U), for
#    every data row.
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 14).Value = $ws.Cells.Item($r, 18).Value2
    $ws.Cells.Item($r, 15).Value = $ws.Cells.Item($r, 19).Value2
    $ws.Cells.Item($r, 16).Value = $ws.Cells.Item($r, 20).Value2
    $ws.Cells.Item($r, 17).Value = $ws.Cells.Item($r, 21).Value2
}

# 4) Rename the shifted former "initial01".."initial04" headers (now at
#    R1:U1) to the "total initial biomass" variant. Their values are left
#    untouched.
$ws.Range("R1").Value = "initial01totmass"
$ws.Range("S1").Value = "initial02totmass"
$ws.Range("T1").Value = "initial03totmass"
$ws.Range("U1").Value = "initial04totmass"

# 5) The freshly inserted N:Q columns don't inherit a bestFit width; nudge
#    them close to the (bestFit) width used by the other "wide" numeric
#    blocks such as the old overwinter block.
$ws.Range("N1:Q1").EntireColumn.ColumnWidth = 10.3

# 6) Reflect the new selection recorded in the sheet view.
$ws.Range("J28").Select()
